$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 188.16667
$ws.Range("J2").Value = 132.25
$ws.Range("L2").Value = 132.25
$ws.Range("N2").Value = -358.25

$ws.Range("H80").Value = 13842.375
$ws.Range("I80").Value = 1755.6
$ws.Range("K80").Value = 5266.799999999999
$ws.Range("M80").Value = -4268.799999999999

$ws.Range("H83").Value = 13842.375
$ws.Range("I83").Value = 1755.6
$ws.Range("K83").Value = 15800.4
$ws.Range("M83").Value = -10808.4

$ws.Range("H116").Value = 12917.714
$ws.Range("J116").Value = 6463.4546
$ws.Range("L116").Value = 6463.4546
$ws.Range("N116").Value = -13347.4546

$ws.Range("H132").Value = 1644.9565
$ws.Range("I132").Value = 1644.9565
$ws.Range("K132").Value = 4934.8695
$ws.Range("M132").Value = -2404.8695

$ws.Range("H137").Value = 4047
$ws.Range("I137").Value = 1400.5
$ws.Range("J137").Value = 5223.222
$ws.Range("K137").Value = 4201.5
$ws.Range("L137").Value = 15669.666
$ws.Range("M137").Value = -1651.5
$ws.Range("N137").Value = -20769.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3672.9644
$ws.Range("I32").Value = 2063.689
$ws.Range("J32").Value = 10256.363
$ws.Range("K32").Value = 2063.689
$ws.Range("L32").Value = 10256.363
$ws.Range("M32").Value = -1776.689
$ws.Range("N32").Value = -10830.363

$ws.Range("H45").Value = 1466.0526
$ws.Range("I45").Value = 1085.8889
$ws.Range("K45").Value = 1085.8889
$ws.Range("M45").Value = -708.8888999999999

$ws.Range("H61").Value = 4533.9565
$ws.Range("I61").Value = 2537.125
$ws.Range("K61").Value = 2537.125
$ws.Range("M61").Value = -2325.125

$ws.Range("H74").Value = 896.14813
$ws.Range("I74").Value = 863.88
$ws.Range("K74").Value = 863.88
$ws.Range("M74").Value = 10.12

$ws.Range("H77").Value = 896.14813
$ws.Range("I77").Value = 863.88
$ws.Range("K77").Value = 4319.4
$ws.Range("M77").Value = 48.60000000000036

$ws.Range("H96").Value = 10344
$ws.Range("J96").Value = 10344
$ws.Range("L96").Value = 10344
$ws.Range("N96").Value = -15836

$ws.Range("H110").Value = 3504.3333
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H132").Value = 1394.4615
$ws.Range("I132").Value = 1054.6957
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 3164.0871
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -634.0870999999997
$ws.Range("N132").Value = -17057.9999

$ws.Range("H136").Value = 4533.9565
$ws.Range("I136").Value = 2537.125
$ws.Range("K136").Value = 7611.375
$ws.Range("M136").Value = -5061.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1228.3077
$ws.Range("I99").Value = 954.8570999999999
$ws.Range("K99").Value = 954.8570999999999
$ws.Range("M99").Value = 543.1429000000001

$ws.Range("H105").Value = 2258
$ws.Range("I105").Value = 2258
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2258
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -511
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 2471.125
$ws.Range("I107").Value = 1294.8334
$ws.Range("K107").Value = 1294.8334
$ws.Range("M107").Value = 625.1666

$ws.Range("H108").Value = 94977.5
$ws.Range("J108").Value = 94977.5
$ws.Range("L108").Value = 94977.5
$ws.Range("N108").Value = -102657.5

$ws.Range("H134").Value = 10150.464
$ws.Range("I134").Value = 10714.714
$ws.Range("K134").Value = 32144.142
$ws.Range("M134").Value = -29609.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 443.66666
$ws.Range("I107").Value = 387.91666
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 387.91666
$ws.Range("L107").Value = 666.6667
$ws.Range("M107").Value = 1532.08334
$ws.Range("N107").Value = -4506.6667

$ws.Range("H134").Value = 828.34485
$ws.Range("I134").Value = 815.1111
$ws.Range("K134").Value = 2445.3333
$ws.Range("M134").Value = 89.66670000000022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 717.1667
$ws.Range("I5").Value = 680.6
$ws.Range("K5").Value = 2041.8
$ws.Range("M5").Value = -1929.8

$ws.Range("H26").Value = 461
$ws.Range("J26").Value = 501.25
$ws.Range("L26").Value = 1503.75
$ws.Range("N26").Value = -2079.75

$ws.Range("H131").Value = 786.39
$ws.Range("J131").Value = 798.38947
$ws.Range("L131").Value = 2395.16841
$ws.Range("N131").Value = -12475.16841

$ws.Range("H135").Value = 717.1667
$ws.Range("I135").Value = 680.6
$ws.Range("K135").Value = 6125.400000000001
$ws.Range("M135").Value = -3590.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2874.4
$ws.Range("I80").Value = 2211
$ws.Range("J80").Value = 3316.6667
$ws.Range("K80").Value = 2211
$ws.Range("L80").Value = 3316.6667
$ws.Range("M80").Value = -1213
$ws.Range("N80").Value = -5312.6667

$ws.Range("H83").Value = 2874.4
$ws.Range("I83").Value = 2211
$ws.Range("J83").Value = 3316.6667
$ws.Range("K83").Value = 11055
$ws.Range("L83").Value = 16583.3335
$ws.Range("M83").Value = -6063
$ws.Range("N83").Value = -26567.3335

$ws.Range("H122").Value = 2332.3333
$ws.Range("J122").Value = 2332.3333
$ws.Range("L122").Value = 6996.999899999999
$ws.Range("N122").Value = -11896.9999

$ws.Range("H132").Value = 1541581.9
$ws.Range("I132").Value = 2264129.8
$ws.Range("K132").Value = 6792389.399999999
$ws.Range("M132").Value = -6789859.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 7233.6924
$ws.Range("I4").Value = 5070.5
$ws.Range("J4").Value = 14444.333
$ws.Range("K4").Value = 5070.5
$ws.Range("L4").Value = 14444.333
$ws.Range("M4").Value = -4957.5
$ws.Range("N4").Value = -14670.333

$ws.Range("H21").Value = 9194.25
$ws.Range("I21").Value = 9000
$ws.Range("J21").Value = 9259
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 9259
$ws.Range("M21").Value = -8826
$ws.Range("N21").Value = -9607

$ws.Range("H28").Value = 7233.6924
$ws.Range("I28").Value = 5070.5
$ws.Range("J28").Value = 14444.333
$ws.Range("K28").Value = 5070.5
$ws.Range("L28").Value = 14444.333
$ws.Range("M28").Value = -4838.5
$ws.Range("N28").Value = -14908.333

$ws.Range("H30").Value = 508
$ws.Range("I30").Value = 508
$ws.Range("K30").Value = 508
$ws.Range("M30").Value = -400

$ws.Range("H37").Value = 7233.6924
$ws.Range("I37").Value = 5070.5
$ws.Range("J37").Value = 14444.333
$ws.Range("K37").Value = 5070.5
$ws.Range("L37").Value = 14444.333
$ws.Range("M37").Value = -4963.5
$ws.Range("N37").Value = -14658.333

$ws.Range("H98").Value = 99999
$ws.Range("J98").Value = 99999
$ws.Range("L98").Value = 99999
$ws.Range("N98").Value = -105989

$ws.Range("H134").Value = 58900
$ws.Range("J134").Value = 58900
$ws.Range("L134").Value = 58900
$ws.Range("N134").Value = -69040

$ws.Range("H136").Value = 3301.625
$ws.Range("I136").Value = 1802.6
$ws.Range("J136").Value = 5800
$ws.Range("K136").Value = 5407.799999999999
$ws.Range("L136").Value = 17400
$ws.Range("M136").Value = -2857.799999999999
$ws.Range("N136").Value = -22500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 28745
$ws.Range("J101").Value = 28745
$ws.Range("L101").Value = 28745
$ws.Range("N101").Value = -35235

$ws.Range("H113").Value = 1111.6923
$ws.Range("I113").Value = 887.75
$ws.Range("K113").Value = 2663.25
$ws.Range("M113").Value = -493.25

$ws.Range("H136").Value = 18520022
$ws.Range("I136").Value = 25253524
$ws.Range("J136").Value = 2893.625
$ws.Range("K136").Value = 75760572
$ws.Range("M136").Value = -75758022
$ws.Range("N136").Value = -13780.875
